# The commit adds an OLE_LINK1 bookmark spanning (almost) the whole body of
# the document, and relocates the auto-generated "_GoBack" bookmark (the
# marker Word drops at the location of the user's last edit) from the very
# end of the document to a point in the middle of a run further up, splitting
# that run's text in the process.

$d = $word.ActiveDocument

# --- 1. Remember where the existing "_GoBack" bookmark currently sits -----
# (it marks the end of the last edit position, right before the final
# paragraph mark) -- this position becomes the end of the new OLE_LINK1
# bookmark.
$goBack = $d.Bookmarks("_GoBack")
$goBackEnd = $goBack.End

# --- 2. Add the OLE_LINK1 bookmark, spanning from the very start of the ---
# document up to that same position.
$oleRange = $d.Range(0, $goBackEnd)
$d.Bookmarks.Add("OLE_LINK1", $oleRange)

# --- 3. Locate the last occurrence of "ng th" before the ProductService ---
# sentence -- this is where the cursor was when the "_GoBack" marker needs
# to move to, splitting the run "ng th" into "ng" and " th".
$scanRange = $d.Range(0, $d.Content.End)
$matchStart = -1
$matchEnd = -1
while ($scanRange.Find.Execute("ng th", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchStart = $scanRange.Start
    $matchEnd = $scanRange.End
    $scanRange.Start = $scanRange.End
    $scanRange.End = $d.Content.End
}
$splitPos = $matchStart + 2

# --- 4. Move "_GoBack": delete the old one and re-add it, collapsed, at ---
# the split point found above.
$goBack2 = $d.Bookmarks("_GoBack")
$goBack2.Delete()

$newGoBackRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $newGoBackRange)
